# Stock update during the sales process: append the newly recorded sales
# (rows 13-18) to the sell_data sheet, and refresh the re-saved E12
# timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E12 picked up a (numerically identical) re-saved fractional value.
$ws.Cells.Item(12, 5).Value = 45815.82793381945

# New sale rows appended after row 12 (id, items, quantities, total, date).
$rows = @(
    @(12, "5AYB",     "3",    36000,  45818.69865769676),
    @(13, "5AYB-1AP", "2-10", 109000, 45818.6991234375),
    @(14, "1AYB",     "2",    7000,   45818.70166679398),
    @(15, "5AYB-1AP", "3-10", 121000, 45818.70304361111),
    @(16, "1M-1L",    "5-4",  97000,  45818.70418662037),
    @(17, "2M",       "5",    140000, 45818.70474743722)
)

$r = 13
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    # "quantities" is always stored as text. A plain digit string (e.g. "3")
    # would otherwise be auto-recognised as a number by the Value setter, so
    # route those through a text formula + paste-as-values round trip to
    # keep them as genuine text without leaving a quote-prefix style behind.
    $qty = $row[2]
    if ($qty -match '^[0-9]+$') {
        $ws.Cells.Item($r, 3).Formula = '="' + $qty + '"'
        $ws.Cells.Item($r, 3).Copy()
        $ws.Cells.Item($r, 3).PasteSpecial(-4163)
    } else {
        $ws.Cells.Item($r, 3).Value = $qty
    }

    $ws.Cells.Item($r, 4).Value = $row[3]

    # Match the existing date-time number format used by column E.
    $ws.Cells.Item($r, 5).NumberFormat = $ws.Cells.Item(2, 5).NumberFormat
    $ws.Cells.Item($r, 5).Value = $row[4]

    $r++
}
